$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 45834
$ws.Range("C25").Value = 69596484
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = "Ronaldo"
$ws.Range("F25").Value = "Estava fazendo triagem Infracommerce"

# Row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 45834
$ws.Range("C26").Value = 69596484
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = "Bruno"
$ws.Range("F26").Value = "Bruno fez sozinho"

# Match formatting of the preceding data rows (date + order-number columns)
$ws.Range("B16").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B26").PasteSpecial(-4122)

